$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2"  = 48.35174455092411
    "B3"  = 64.0552961547386
    "B4"  = 72.30708663865317
    "H5"  = 96.15680975952145
    "H6"  = 96.16345600422915
    "H7"  = 96.1495277258243
    "C8"  = 95.73711576545456
    "C9"  = 92.23797695844402
    "C10" = 94.77336859032364
    "D11" = 99.03271399899288
    "D12" = 98.54814472606226
    "D13" = 99.42057136409449
    "E14" = 98.82842695984328
    "E15" = 98.95398693150035
    "E16" = 98.89990867616731
    "F17" = 98.43808881858786
    "F18" = 98.49999765928868
    "F19" = 98.43494036128116
    "G20" = 97.51668404291152
    "G21" = 97.59771320191248
    "G22" = 97.55445373862165
    "B23" = 75.48004170624002
    "B24" = 76.54837578481774
    "H25" = 96.17984708002635
    "H26" = 96.15440184946873
    "C27" = 94.73498098548896
    "C28" = 94.41306125100357
    "D29" = 98.89751425347822
    "D30" = 99.11369857655407
    "E31" = 98.97058365506778
    "E32" = 98.93755316944949
    "F33" = 98.48205786916955
    "F34" = 98.55505160818977
    "G35" = 97.49762375607466
    "G36" = 97.51288793677189
    "B37" = 70.70973037849292
    "B38" = 88.81536273227179
    "H39" = 96.1496483056405
    "H40" = 96.07435739501985
    "C41" = 95.10073366532869
    "C42" = 93.20838616041503
    "D43" = 98.91440946149281
    "D44" = 99.05512027426806
    "E45" = 98.87523201920094
    "E46" = 98.34790327889318
    "F47" = 98.54440849522305
    "F48" = 98.39651592086329
    "G49" = 97.50137930780609
    "G50" = 97.62363985457112
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
